$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# The sheet is a correlation matrix. A new variable
# "2008-9RecessionDummyVar" is being inserted as the second-to-last
# row/column (just before "AnnualizedMoM-CPI-Inflation", which shifts
# from column Q / row 17 to column R / row 18).
# ------------------------------------------------------------------

# New correlation values between "2008-9RecessionDummyVar" and each of
# the existing 15 variables (rows 2-16 / columns B-P), in order.
$newVals = @(
    -0.3040967956544474,
    -0.2014285535872592,
    -0.2400181772862542,
    -0.01743733957108285,
    -0.01769392986876932,
    0.04402117411545403,
    -0.02468419043730151,
    -0.03679216369613469,
    -0.2677739025551792,
    0.03518189690949631,
    -0.06424763057618306,
    0.04909288272459091,
    -0.1237936718063458,
    0.1013215804987397,
    0.03525849430657787
)

# Correlation between the new dummy var and "AnnualizedMoM-CPI-Inflation"
$corrWithCpi = -0.2932899553554737

# ------------------------------------------------------------------
# Step 1: shift existing column Q (B..Q data, header in row1, rows1-17)
# out to column R, preserving formatting, so we can reuse column Q for
# the new variable.
# ------------------------------------------------------------------
$ws.Range("Q1:Q17").Copy($ws.Range("R1:R17"))

# ------------------------------------------------------------------
# Step 2: shift existing row 17 (A..Q, now A..R since the column copy
# above already populated R17) out to row 18, preserving formatting.
# ------------------------------------------------------------------
$ws.Range("A17:R17").Copy($ws.Range("A18:R18"))

# ------------------------------------------------------------------
# Step 3: populate the new column Q (rows 1-16) with the new variable's
# header + correlation values against the other (non-dummy) variables.
# ------------------------------------------------------------------
$ws.Range("Q1").Value = "2008-9RecessionDummyVar"

for ($i = 0; $i -lt $newVals.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 17).Value = $newVals[$i]
}

# ------------------------------------------------------------------
# Step 4: populate the new row 17 (A..P) with the new variable's label
# + the same correlation values (symmetric matrix), plus the
# self-correlation (Q17 = 1) and the correlation against
# "AnnualizedMoM-CPI-Inflation" (R17).
# ------------------------------------------------------------------
$ws.Range("A2").Copy($ws.Range("A17"))
$ws.Range("A17").Value = "2008-9RecessionDummyVar"

for ($i = 0; $i -lt $newVals.Length; $i++) {
    $col = $i + 2
    $ws.Cells.Item(17, $col).Value = $newVals[$i]
}

$ws.Cells.Item(17, 17).Value = 1
$ws.Cells.Item(17, 18).Value = $corrWithCpi

# ------------------------------------------------------------------
# Step 5: fix up the shifted row/column intersection cells that the
# block-copies in steps 1-2 could not have populated correctly on
# their own (Q18, the correlation between CPI-inflation and the new
# dummy var, and R18 which must remain the CPI self-correlation = 1).
# ------------------------------------------------------------------
$ws.Cells.Item(18, 17).Value = $corrWithCpi
$ws.Cells.Item(18, 18).Value = 1

$wb.Save()
